$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Height")

# Insert a new row above the current row 13 (MATO), shifting existing
# rows 13-40 down to 14-41, then populate the new row with the
# EWV/LN entry.
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = "EWV/LN"
$ws.Cells.Item(13, 2).Value = "H:1"
$ws.Cells.Item(13, 3).Value = 1
